$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of (row, [C,D,E,F,G] text values) — new dataset with
# corrected column names. Values are written as text (quote-prefixed) so
# Excel preserves the exact formatted strings (trailing zeros, fixed
# decimal places) instead of silently re-parsing them as numbers.
$rows = @(
    [pscustomobject]@{ Row = 2;  Values = @("51.80", "3.20", "1.00",  "0.0860", "0.1309") }
    [pscustomobject]@{ Row = 3;  Values = @("11.40", "4.20", "74.20", "0.7505", "0.7520") }
    [pscustomobject]@{ Row = 4;  Values = @("53.60", "8.20", "60.50", "0.9420", "0.9262") }
    [pscustomobject]@{ Row = 5;  Values = @("51.80", "3.20", "1.00",  "0.2630", "0.4089") }
    [pscustomobject]@{ Row = 6;  Values = @("11.50", "4.00", "53.20", "1.6776", "1.6784") }
    [pscustomobject]@{ Row = 7;  Values = @("54.70", "9.50", "41.40", "2.2000", "2.1423") }
    [pscustomobject]@{ Row = 8;  Values = @("51.80", "3.20", "1.00",  "0.6520", "0.9920") }
    [pscustomobject]@{ Row = 9;  Values = @("54.80", "4.40", "14.90", "3.1613", "3.1456") }
    [pscustomobject]@{ Row = 10; Values = @("34.10", "9.90", "63.30", "4.1000", "4.0157") }
)

foreach ($entry in $rows) {
    $row = $entry.Row
    $values = $entry.Values
    for ($col = 3; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = "'" + $values[$col - 3]
    }
}
